$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 513277.7
$ws.Range("I33").Value = 641274.4
$ws.Range("K33").Value = 641274.4
$ws.Range("M33").Value = -641045.4
$ws.Range("H100").Value = 2971.1025
$ws.Range("I100").Value = 2028.125
$ws.Range("K100").Value = 2028.125
$ws.Range("M100").Value = -1487.125
$ws.Range("H112").Value = 780260.3
$ws.Range("I112").Value = 2062.25
$ws.Range("J112").Value = 1402818.8
$ws.Range("K112").Value = 6186.75
$ws.Range("L112").Value = 4208456.4
$ws.Range("M112").Value = -5078.75
$ws.Range("N112").Value = -4210672.4
$ws.Range("H121").Value = 2177.652
$ws.Range("J121").Value = 2203.9092
$ws.Range("L121").Value = 6611.7276
$ws.Range("N121").Value = -10105.7276
$ws.Range("H135").Value = 1104.6111
$ws.Range("I135").Value = 999
$ws.Range("J135").Value = 2900
$ws.Range("K135").Value = 8991
$ws.Range("L135").Value = 26100
$ws.Range("M135").Value = -6456
$ws.Range("N135").Value = -31170
$ws.Range("H138").Value = 4227.95
$ws.Range("I138").Value = 2026.8182
$ws.Range("K138").Value = 6080.4546
$ws.Range("M138").Value = -940.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17641014
$ws.Range("I32").Value = 18586952
$ws.Range("K32").Value = 18586952
$ws.Range("M32").Value = -18586665
$ws.Range("H63").Value = 3582.15
$ws.Range("I63").Value = 1972.0769
$ws.Range("J63").Value = 6572.2856
$ws.Range("K63").Value = 1972.0769
$ws.Range("L63").Value = 6572.2856
$ws.Range("M63").Value = -1286.0769
$ws.Range("N63").Value = -7944.2856
$ws.Range("H66").Value = 3582.15
$ws.Range("I66").Value = 1972.0769
$ws.Range("J66").Value = 6572.2856
$ws.Range("K66").Value = 9860.3845
$ws.Range("L66").Value = 32861.428
$ws.Range("M66").Value = -6428.3845
$ws.Range("N66").Value = -39725.428
$ws.Range("H74").Value = 3450
$ws.Range("I74").Value = 3388.889
$ws.Range("K74").Value = 3388.889
$ws.Range("M74").Value = -2514.889
$ws.Range("H77").Value = 3450
$ws.Range("I77").Value = 3388.889
$ws.Range("K77").Value = 16944.445
$ws.Range("M77").Value = -12576.445
$ws.Range("H132").Value = 4592.8667
$ws.Range("I132").Value = 4428.0713
$ws.Range("J132").Value = 6900
$ws.Range("K132").Value = 13284.2139
$ws.Range("L132").Value = 20700
$ws.Range("M132").Value = -10754.2139
$ws.Range("N132").Value = -25760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 625
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -1346
$ws.Range("H99").Value = 1768.2069
$ws.Range("I99").Value = 1773.2084
$ws.Range("J99").Value = 1744.2
$ws.Range("K99").Value = 1773.2084
$ws.Range("L99").Value = 1744.2
$ws.Range("M99").Value = -275.2084
$ws.Range("N99").Value = -4740.2
$ws.Range("H107").Value = 3465.6924
$ws.Range("I107").Value = 3491.8262
$ws.Range("J107").Value = 3265.3333
$ws.Range("K107").Value = 3491.8262
$ws.Range("L107").Value = 3265.3333
$ws.Range("M107").Value = -1571.8262
$ws.Range("N107").Value = -7105.3333
$ws.Range("H134").Value = 5131764.5
$ws.Range("I134").Value = 6064202.5
$ws.Range("J134").Value = 3356.5
$ws.Range("K134").Value = 18192607.5
$ws.Range("L134").Value = 10069.5
$ws.Range("M134").Value = -18190072.5
$ws.Range("N134").Value = -15139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3078.4143
$ws.Range("I31").Value = 1632.2778
$ws.Range("J31").Value = 4609.6177
$ws.Range("K31").Value = 1632.2778
$ws.Range("L31").Value = 4609.6177
$ws.Range("M31").Value = -1337.2778
$ws.Range("N31").Value = -5199.6177
$ws.Range("H34").Value = 3078.4143
$ws.Range("I34").Value = 1632.2778
$ws.Range("J34").Value = 4609.6177
$ws.Range("K34").Value = 1632.2778
$ws.Range("L34").Value = 4609.6177
$ws.Range("M34").Value = -1430.2778
$ws.Range("N34").Value = -5013.6177
$ws.Range("H58").Value = 3267.111
$ws.Range("I58").Value = 2679.818
$ws.Range("K58").Value = 2679.818
$ws.Range("M58").Value = -2476.818
$ws.Range("H92").Value = 62663.332
$ws.Range("J92").Value = 62663.332
$ws.Range("L92").Value = 62663.332
$ws.Range("N92").Value = -67655.33199999999
$ws.Range("H96").Value = 214553.2
$ws.Range("J96").Value = 214553.2
$ws.Range("L96").Value = 214553.2
$ws.Range("N96").Value = -220045.2
$ws.Range("H132").Value = 1618
$ws.Range("I132").Value = 1100.5333
$ws.Range("K132").Value = 3301.5999
$ws.Range("M132").Value = -771.5999000000002
$ws.Range("H134").Value = 2150.5386
$ws.Range("I134").Value = 2127.1904
$ws.Range("K134").Value = 6381.5712
$ws.Range("M134").Value = -3846.5712
$ws.Range("H136").Value = 3267.111
$ws.Range("I136").Value = 2679.818
$ws.Range("K136").Value = 8039.454000000001
$ws.Range("M136").Value = -5489.454000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1166.3334
$ws.Range("I5").Value = 999
$ws.Range("J5").Value = 1199.8
$ws.Range("K5").Value = 2997
$ws.Range("L5").Value = 3599.4
$ws.Range("M5").Value = -2885
$ws.Range("N5").Value = -3823.4
$ws.Range("H12").Value = 555617.5
$ws.Range("I12").Value = 33
$ws.Range("J12").Value = 833409.75
$ws.Range("K12").Value = 99
$ws.Range("L12").Value = 2500229.25
$ws.Range("M12").Value = 74
$ws.Range("N12").Value = -2500575.25
$ws.Range("H107").Value = 881.2727
$ws.Range("I107").Value = 769.5
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 2308.5
$ws.Range("L107").Value = 5997
$ws.Range("M107").Value = -388.5
$ws.Range("N107").Value = -9837
$ws.Range("H113").Value = 1490.7693
$ws.Range("I113").Value = 1138.6
$ws.Range("J113").Value = 1710.875
$ws.Range("K113").Value = 3415.8
$ws.Range("L113").Value = 5132.625
$ws.Range("M113").Value = -1245.8
$ws.Range("N113").Value = -9472.625
$ws.Range("H122").Value = 634198.25
$ws.Range("I122").Value = 798
$ws.Range("J122").Value = 996141.3
$ws.Range("K122").Value = 7182
$ws.Range("L122").Value = 8965271.700000001
$ws.Range("M122").Value = -4732
$ws.Range("N122").Value = -8970171.700000001
$ws.Range("H135").Value = 1166.3334
$ws.Range("I135").Value = 999
$ws.Range("J135").Value = 1199.8
$ws.Range("K135").Value = 8991
$ws.Range("L135").Value = 10798.2
$ws.Range("M135").Value = -6456
$ws.Range("N135").Value = -15868.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1739.2
$ws.Range("J107").Value = 1549.25
$ws.Range("L107").Value = 1549.25
$ws.Range("N107").Value = -5389.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40540
$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -41872
$ws.Range("H132").Value = 5740
$ws.Range("I132").Value = 5479.8335
$ws.Range("J132").Value = 6260.3335
$ws.Range("K132").Value = 16439.5005
$ws.Range("L132").Value = 18781.0005
$ws.Range("M132").Value = -13909.5005
$ws.Range("N132").Value = -23841.0005
$ws.Range("H134").Value = 99995.664
$ws.Range("J134").Value = 99995.664
$ws.Range("L134").Value = 99995.664
$ws.Range("N134").Value = -110135.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 55436.555
$ws.Range("J119").Value = 55436.555
$ws.Range("L119").Value = 55436.555
$ws.Range("N119").Value = -65112.555
$ws.Range("H136").Value = 66538.875
$ws.Range("I136").Value = 3132.5
$ws.Range("J136").Value = 172216.17
$ws.Range("K136").Value = 9397.5
$ws.Range("L136").Value = 516648.51
$ws.Range("M136").Value = -6847.5
$ws.Range("N136").Value = -521748.51

Write-Host "Applied all updates"